# Maj Comparaison Frame Invictus-Optimus
# Updates the "Invictus (kg)" column (D) on the Frame_Comparaison sheet with
# the latest mass figures, lets the dependent Gain columns (E/F) and the
# Total row recompute, and leaves the cursor on the last-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Front Hoop Bracing Support (row 9)
$ws.Range("D9").Formula = "=2*(0.34+0.632+0.123+0.184+0.102+0.19+0.286+0.285)"

# Tubes Non Structurant Cellule Avant (hors plancher) (row 10)
$ws.Range("D10").Value = 0

# Front Hoop (row 11)
$ws.Range("D11").Value = 1.806

# Plancher Cellule Avant (row 12)
$ws.Range("D12").Formula = "=0.236+0.237+2*0.126"

# Side Impact Structure (row 13)
$ws.Range("D13").Formula = "=2*(0.792+0.731+0.275+0.615)"

# Plancher Cockpit (row 14)
$ws.Range("D14").Formula = "=0.575*2"

# Main Hoop (row 16)
$ws.Range("D16").Value = 3.5

# Leave the selection where the author ended up
$ws.Range("D17").Select() | Out-Null
